# Refresh Betfair Back/Lay odds for Jogos_do_Dia_Betfair_Back_Lay_2026-02-18.xlsx
# Updates numeric odds values (columns F:AO) for rows 2-11 on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 2.72
$ws.Range("G2").Value = 3.4
$ws.Range("I2").Value = 3.15
$ws.Range("J2").Value = 3.2
$ws.Range("K2").Value = 3.8
$ws.Range("L2").Value = 1.01
$ws.Range("M2").Value = 1.07
$ws.Range("N2").Value = 3.1
$ws.Range("O2").Value = 1.32
$ws.Range("P2").Value = 1.84
$ws.Range("Q2").Value = 1.94
$ws.Range("R2").Value = 1.32
$ws.Range("S2").Value = 3.05
$ws.Range("T2").Value = 1.62
$ws.Range("U2").Value = 1.94
$ws.Range("V2").Value = 1.47
$ws.Range("W2").Value = 1.42
$ws.Range("X2").Value = 17
$ws.Range("Y2").Value = 13.5
$ws.Range("Z2").Value = 22
$ws.Range("AA2").Value = 50
$ws.Range("AB2").Value = 14
$ws.Range("AC2").Value = 9.4
$ws.Range("AD2").Value = 15
$ws.Range("AE2").Value = 38
$ws.Range("AF2").Value = 24
$ws.Range("AG2").Value = 15.5
$ws.Range("AH2").Value = 21
$ws.Range("AI2").Value = 55
$ws.Range("AJ2").Value = 60
$ws.Range("AK2").Value = 40
$ws.Range("AL2").Value = 55
$ws.Range("AM2").Value = 120
$ws.Range("AN2").Value = 36
$ws.Range("AO2").Value = 32

# Row 3
$ws.Range("P3").Value = 1.81

# Row 4
$ws.Range("F4").Value = 2.36
$ws.Range("G4").Value = 2.94
$ws.Range("H4").Value = 2.82
$ws.Range("I4").Value = 3.6
$ws.Range("J4").Value = 3.1
$ws.Range("K4").Value = 4
$ws.Range("P4").Value = 2
$ws.Range("Q4").Value = 1.78

# Row 5
$ws.Range("G5").Value = 9.2
$ws.Range("I5").Value = 1.43
$ws.Range("J5").Value = 5.3
$ws.Range("O5").Value = 1.21
$ws.Range("P5").Value = 2.24
$ws.Range("Q5").Value = 1.67
$ws.Range("S5").Value = 2.68
$ws.Range("U5").Value = 2
$ws.Range("AB5").Value = 32
$ws.Range("AH5").Value = 26
$ws.Range("AJ5").Value = 1000
$ws.Range("AK5").Value = 1000
$ws.Range("AL5").Value = 1000
$ws.Range("AM5").Value = 1000
$ws.Range("AN5").Value = 1000
$ws.Range("AO5").Value = 6

# Row 6
$ws.Range("F6").Value = 4.2
$ws.Range("M6").Value = 1.04
$ws.Range("O6").Value = 1.24
$ws.Range("P6").Value = 2.32
$ws.Range("Q6").Value = 1.71
$ws.Range("R6").Value = 1.52
$ws.Range("S6").Value = 2.78
$ws.Range("U6").Value = 2.38
$ws.Range("X6").Value = 23
$ws.Range("AB6").Value = 1000
$ws.Range("AC6").Value = 9.6
$ws.Range("AF6").Value = 1000
$ws.Range("AJ6").Value = 1000
$ws.Range("AL6").Value = 1000
$ws.Range("AN6").Value = 44
$ws.Range("AO6").Value = 9.8

# Row 8
$ws.Range("N8").Value = 4.7
$ws.Range("O8").Value = 1.24
$ws.Range("P8").Value = 2.26
$ws.Range("R8").Value = 1.49
$ws.Range("T8").Value = 2.28
$ws.Range("X8").Value = 21
$ws.Range("AC8").Value = 13.5
$ws.Range("AF8").Value = 140
$ws.Range("AG8").Value = 50
$ws.Range("AH8").Value = 36
$ws.Range("AI8").Value = 48
$ws.Range("AK8").Value = 290
$ws.Range("AL8").Value = 210
$ws.Range("AM8").Value = 410
$ws.Range("AO8").Value = 5.3

# Row 9
$ws.Range("F9").Value = 2.68
$ws.Range("G9").Value = 2.72
$ws.Range("I9").Value = 2.92
$ws.Range("N9").Value = 3.85
$ws.Range("S9").Value = 3.45
$ws.Range("Z9").Value = 20
$ws.Range("AL9").Value = 40

# Row 10
$ws.Range("G10").Value = 4.2
$ws.Range("H10").Value = 1.97
$ws.Range("I10").Value = 1.98
$ws.Range("N10").Value = 4.9
$ws.Range("O10").Value = 1.2
$ws.Range("P10").Value = 2.36
$ws.Range("Q10").Value = 1.69
$ws.Range("R10").Value = 1.49
$ws.Range("S10").Value = 2.6
$ws.Range("T10").Value = 1.6
$ws.Range("U10").Value = 2.28
$ws.Range("AA10").Value = 24
$ws.Range("AB10").Value = 23
$ws.Range("AC10").Value = 10
$ws.Range("AE10").Value = 21
$ws.Range("AF10").Value = 38
$ws.Range("AG10").Value = 19.5
$ws.Range("AK10").Value = 50
$ws.Range("AM10").Value = 70

# Row 11
$ws.Range("F11").Value = 5.1
$ws.Range("G11").Value = 5.3
$ws.Range("H11").Value = 1.72
$ws.Range("I11").Value = 1.74
$ws.Range("J11").Value = 4.4
$ws.Range("N11").Value = 5.9
$ws.Range("P11").Value = 2.66
$ws.Range("R11").Value = 1.68
$ws.Range("S11").Value = 2.34
$ws.Range("U11").Value = 2.5
$ws.Range("AF11").Value = 55
$ws.Range("AL11").Value = 50
$ws.Range("AM11").Value = 75
$ws.Range("AN11").Value = 42
$ws.Range("AO11").Value = 6.8
